# Applies the "Updated cryptos list" data refresh described by the commit diff.
# Most cells only change the Price (D) and/or Volume(1h) (E) text; rows 47/48
# additionally swap their Coin/Link/Price/Volume content (dogwifhat <-> OKB).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Range("D2").Value = "75.055.34"
$ws.Range("E2").Value = "  +1.01%  "

# Row 3
$ws.Range("D3").Value = "2.853.14"
$ws.Range("E3").Value = "  +10.49%  "

# Row 4
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").Value = "'604.48"
$ws.Range("E5").Value = "  +3.87%  "

# Row 6
$ws.Range("D6").Value = "'188.81"
$ws.Range("E6").Value = "  +2.60%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").Value = "'0.560"
$ws.Range("E8").Value = "  +4.93%  "

# Row 9
$ws.Range("E9").Value = "  -7.55%  "

# Row 10
$ws.Range("D10").Value = "2.849.74"
$ws.Range("E10").Value = "  +10.39%  "

# Row 11
$ws.Range("E11").Value = "  +0.22%  "

# Row 12
$ws.Range("D12").Value = "'0.373"
$ws.Range("E12").Value = "  +4.55%  "

# Row 13
$ws.Range("D13").Value = "'4.95"
$ws.Range("E13").Value = "  +3.03%  "

# Row 14
$ws.Range("D14").Value = "3.380.58"
$ws.Range("E14").Value = "  +10.55%  "

# Row 15
$ws.Range("D15").Value = "74.978.00"
$ws.Range("E15").Value = "  +1.08%  "

# Row 16
$ws.Range("D16").Value = "'27.70"
$ws.Range("E16").Value = "  +6.04%  "

# Row 17
$ws.Range("D17").Value = "'0.0000190"
$ws.Range("E17").Value = "  -2.07%  "

# Row 18
$ws.Range("D18").Value = "2.855.49"
$ws.Range("E18").Value = "  +11.28%  "

# Row 19
$ws.Range("D19").Value = "'9.21"
$ws.Range("E19").Value = "  +13.25%  "

# Row 20
$ws.Range("D20").Value = "'12.57"
$ws.Range("E20").Value = "  +7.48%  "

# Row 21
$ws.Range("D21").Value = "'379.20"
$ws.Range("E21").Value = "  +4.30%  "

# Row 22
$ws.Range("D22").Value = "'2.30"
$ws.Range("E22").Value = "  +1.10%  "

# Row 23
$ws.Range("D23").Value = "'4.17"
$ws.Range("E23").Value = "  +2.34%  "

# Row 24
$ws.Range("D24").Value = "'6.22"
$ws.Range("E24").Value = "  -0.07%  "

# Row 25
$ws.Range("D25").Value = "'71.26"
$ws.Range("E25").Value = "  +3.43%  "

# Row 26
$ws.Range("E26").Value = "  +0.05%  "

# Row 27
$ws.Range("D27").Value = "'4.27"
$ws.Range("E27").Value = "  +2.62%  "

# Row 28
$ws.Range("D28").Value = "3.009.79"
$ws.Range("E28").Value = "  +11.04%  "

# Row 29
$ws.Range("D29").Value = "'9.71"
$ws.Range("E29").Value = "  +5.80%  "

# Row 30
$ws.Range("D30").Value = "'0.0000105"
$ws.Range("E30").Value = "  +12.08%  "

# Row 31
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  +0.11%  "

# Row 32
$ws.Range("D32").Value = "'533.21"
$ws.Range("E32").Value = "  +6.05%  "

# Row 33
$ws.Range("D33").Value = "'1.42"
$ws.Range("E33").Value = "  +6.68%  "

# Row 34
$ws.Range("D34").Value = "'8.00"
$ws.Range("E34").Value = "  +1.94%  "

# Row 35
$ws.Range("D35").Value = "'1.84"
$ws.Range("E35").Value = "  +7.95%  "

# Row 36
$ws.Range("E36").Value = "  -0.06%  "

# Row 37
$ws.Range("E37").Value = "  +2.39%  "

# Row 38
$ws.Range("D38").Value = "'20.38"
$ws.Range("E38").Value = "  +6.67%  "

# Row 39
$ws.Range("D39").Value = "'162.59"
$ws.Range("E39").Value = "  +1.69%  "

# Row 40
$ws.Range("D40").Value = "'19.31"
$ws.Range("E40").Value = "  -0.28%  "

# Row 41
$ws.Range("D41").Value = "'184.04"
$ws.Range("E41").Value = "  +24.00%  "

# Row 42
$ws.Range("E42").Value = "  +0.08%  "

# Row 43
$ws.Range("D43").Value = "'5.17"
$ws.Range("E43").Value = "  +6.09%  "

# Row 44
$ws.Range("D44").Value = "'0.344"
$ws.Range("E44").Value = "  +8.37%  "

# Row 45
$ws.Range("D45").Value = "'1.70"
$ws.Range("E45").Value = "  +1.49%  "

# Row 46
$ws.Range("D46").Value = "'1.27"
$ws.Range("E46").Value = "  +10.97%  "

# Row 47
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").Value = "'40.00"
$ws.Range("E47").Value = "  +1.69%  "

# Row 48
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").Value = "'2.40"
$ws.Range("E48").Value = "  -0.39%  "

# Row 49
$ws.Range("D49").Value = "'0.0864"
$ws.Range("E49").Value = "  +7.97%  "

# Row 50
$ws.Range("D50").Value = "'0.578"
$ws.Range("E50").Value = "  +11.62%  "

# Row 51
$ws.Range("D51").Value = "'3.78"
$ws.Range("E51").Value = "  +5.58%  "
